# Apply the "Second Commit without Extent" edit:
#  - Fix the typo'd password in D2 (jaga@298 -> jaga@2983), in place so the
#    shared-string table keeps uniqueCount=8 (just count bumps to 12).
#  - Duplicate row 2 (Y / SignIn / testjaga002@gmail.com / jaga@2983) into a
#    new row 3, including the Hyperlink styling + live hyperlinks on C3/D3.
#  - Move the saved selection to D2 and grow the sheet dimension to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the password typo in place -------------------------------------
# Setting the value on the cell that already (uniquely) owns the shared
# string rewrites that shared-string entry in place instead of allocating a
# new one, so uniqueCount stays 8 while usage count grows.
$ws.Range("D2").Value = "jaga@2983"

# --- 2. Duplicate row 2 into row 3 ------------------------------------------
# Read back with Value() (method form) - the bare property getter isn't
# reliable for round-tripping text through this bridge.
$rowLabel    = $ws.Range("A2").Value()
$runMode     = $ws.Range("B2").Value()
$userName    = $ws.Range("C2").Value()
$password    = $ws.Range("D2").Value()

$ws.Range("A3").Value = $rowLabel
$ws.Range("B3").Value = $runMode
$ws.Range("C3").Value = $userName
$ws.Range("D3").Value = $password

# --- 3. Re-create the live mailto hyperlinks on the new row -----------------
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:" + $userName)
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:" + $password)

# Adding a hyperlink stamps its own direct-format style; pull the new cells
# back onto the workbook's shared "Hyperlink" cell style so C3/D3 match C2/D2.
$ws.Range("C3:D3").Style = "Hyperlink"

# --- 4. Match Excel's bookkeeping for the edit -------------------------------
# Dimension auto-extends to A1:D3 once row 3 has data; explicitly move the
# remembered selection to D2 (matches the saved sheetView state).
$ws.Range("D2").Select()
